$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a text number format for column D while writing so strings like
# "1.001" or "1.000" are not auto-converted into numeric values by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.485.19"
$ws.Range("E2").Value = "  -2.31%  "

$ws.Range("D3").Value = "1.866.10"
$ws.Range("E3").Value = "  -2.50%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "329.37"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").Value = "0.4714"
$ws.Range("E7").Value = "  +0.98%  "

$ws.Range("D8").Value = "0.3977"
$ws.Range("E8").Value = "  -1.03%  "

$ws.Range("D9").Value = "47.33"
$ws.Range("E9").Value = "  -10.76%  "

$ws.Range("D10").Value = "0.08036"
$ws.Range("E10").Value = "  -4.36%  "

$ws.Range("E11").Value = "  -2.28%  "

$ws.Range("D12").Value = "21.63"
$ws.Range("E12").Value = "  -2.27%  "

$ws.Range("D13").Value = "1.850.90"
$ws.Range("E13").Value = "  -2.20%  "

$ws.Range("D14").Value = "5.959"
$ws.Range("E14").Value = "  -1.69%  "

$ws.Range("D15").Value = "7.195"
$ws.Range("E15").Value = "  -3.18%  "

$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").Value = "86.60"
$ws.Range("E17").Value = "  -3.35%  "

$ws.Range("E18").Value = "  -2.56%  "

$ws.Range("D19").Value = "0.06561"
$ws.Range("E19").Value = "  -0.71%  "

$ws.Range("D20").Value = "17.32"
$ws.Range("E20").Value = "  -3.36%  "

$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").Value = "5.512"
$ws.Range("E22").Value = "  -4.05%  "

$ws.Range("D23").Value = "27.508.28"
$ws.Range("E23").Value = "  -2.15%  "

$ws.Range("D24").Value = "10.97"
$ws.Range("E24").Value = "  -1.93%  "

$ws.Range("D25").Value = "2.300"
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").Value = "2.068.06"
$ws.Range("E26").Value = "  -2.46%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "154.46"
$ws.Range("E27").Value = "  +0.92%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "20.28"
$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("D29").Value = "2.090"
$ws.Range("E29").Value = "  -2.35%  "

$ws.Range("D30").Value = "5.545"
$ws.Range("E30").Value = "  -3.46%  "

$ws.Range("D31").Value = "122.28"
$ws.Range("E31").Value = "  -0.93%  "

$ws.Range("D32").Value = "0.09498"
$ws.Range("E32").Value = "  -1.68%  "

$ws.Range("D33").Value = "0.9591"
$ws.Range("E33").Value = "  -1.83%  "

$ws.Range("D34").Value = "1.471"
$ws.Range("E34").Value = "  +1.98%  "

$ws.Range("D35").Value = "3.592"
$ws.Range("E35").Value = "  -1.51%  "

$ws.Range("D36").Value = "5.301"
$ws.Range("E36").Value = "  -4.37%  "

$ws.Range("D37").Value = "0.06086"
$ws.Range("E37").Value = "  -1.54%  "

$ws.Range("D38").Value = "0.02246"
$ws.Range("E38").Value = "  -2.35%  "

$ws.Range("D39").Value = "1.217"
$ws.Range("E39").Value = "  -4.14%  "

$ws.Range("D40").Value = "8.096"
$ws.Range("E40").Value = "  -8.57%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5985"
$ws.Range("E41").Value = "  -2.99%  "

$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").Value = "0.9995"
$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("D43").Value = "0.1903"
$ws.Range("E43").Value = "  -0.20%  "

$ws.Range("E44").Value = "  -6.41%  "

$ws.Range("D45").Value = "1.265"
$ws.Range("E45").Value = "  -3.18%  "

$ws.Range("D46").Value = "0.5700"
$ws.Range("E46").Value = "  -2.95%  "

$ws.Range("D47").Value = "12.19"
$ws.Range("E47").Value = "  -5.15%  "

$ws.Range("D48").Value = "3.415"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("D49").Value = "1.940"
$ws.Range("E49").Value = "  -4.12%  "

$ws.Range("D50").Value = "0.06779"
$ws.Range("E50").Value = "  -1.84%  "

$ws.Range("D51").Value = "110.35"
$ws.Range("E51").Value = "  -1.17%  "

# Restore column D style to the sheet default (no explicit style index)
# now that the text values are safely stored, so no stray number format
# is left behind on the cells.
$dRange.Style = $ws.Range("A1").Style
